$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "return" description for /user/register to mention user_exists
$ws.Range("E4").Value = "200，user对象, register_success,register_failed,wrong_mail_code,user_exists"

# Add description for /ifusernameexists/:name row
$ws.Range("F6").Value = "查询用户是否存在"

# Add description for /changepwd/:mail/:code/:pwd row
$ws.Range("F7").Value = "修改密码"

# Move the active selection to F7, matching the final saved state
$ws.Range("F7").Select()
